$wb = $excel.ActiveWorkbook

# --- Sheet "Rushing" ---
$ws = $wb.Worksheets.Item("Rushing")

# Row 2 - J.Allen
$ws.Range("C2").Value = 37
$ws.Range("D2").Value = 41
$ws.Range("E2").Value = 47
$ws.Range("F2").Value = 34

# Row 3 - D.Singletary
$ws.Range("C3").Value = 112
$ws.Range("D3").Value = 84
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = 46

# Row 8 - I.McKenzie
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 7
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 6

# --- Sheet "Receiving" ---
$ws2 = $wb.Worksheets.Item("Receiving")

# Row 2 - D.Singletary
$ws2.Range("C2").Value = 56
$ws2.Range("D2").Value = 45
$ws2.Range("G2").Value = 6
$ws2.Range("H2").Value = 6

# Row 5 - R.Gilliam
$ws2.Range("C5").Value = 5
$ws2.Range("D5").Value = 4

# Row 6 - S.Diggs
$ws2.Range("C6").Value = 127
$ws2.Range("D6").Value = 91
$ws2.Range("E6").Value = 36
$ws2.Range("G6").Value = 29

# Row 7 - E.Sanders
$ws2.Range("C7").Value = 48
$ws2.Range("D7").Value = 33

# Row 8 - C.Beasley
$ws2.Range("C8").Value = 111
$ws2.Range("D8").Value = 83
$ws2.Range("E8").Value = 11
$ws2.Range("F8").Value = 6

# Row 9 - G.Davis
$ws2.Range("C9").Value = 45
$ws2.Range("D9").Value = 27
$ws2.Range("E9").Value = 29
$ws2.Range("F9").Value = 17
$ws2.Range("G9").Value = 20
$ws2.Range("H9").Value = 12

# Row 10 - I.McKenzie
$ws2.Range("C10").Value = 16
$ws2.Range("D10").Value = 13

# Row 12 - D.Knox
$ws2.Range("C12").Value = 58
$ws2.Range("D12").Value = 45
